$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data point: row 2 gets a Prediction (G2) value, which makes
# H2 (Residual = Prediction - Actual) and a new I2 (Residual % = H2/G2) compute.
$ws.Range("G2").Value = 7.08
$ws.Range("H2").Formula = "=G2-F2"
$ws.Range("I2").Formula = "=H2/G2"

# New week row (row 3) appended to the results table.
$ws.Range("A3").Value = 2
$ws.Range("B3").Formula = "=B2+7"
$ws.Range("C3").Formula = "=C2+7"
$ws.Range("D3").Value = 5.91940035675141
$ws.Range("E3").Value = 6.27509964324863
$ws.Range("F3").Value = 6.09725
$ws.Range("H3").Formula = "=G3-F3"
$ws.Range("I3").Formula = "=H3/G3"

# Match number formatting (date style) used by the existing date columns
# by copying the formatting from row 2 instead of assigning a fresh format.
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I3").Select()
